$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.093.51'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.652.99'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.76'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5253'
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2669'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06364'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.53'
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07684'
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.597'
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.691.19'
$ws.Range('E13').Value = '  +1.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.880.09'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5613'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8225'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.62'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.092.60'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.684'
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.31'
$ws.Range('E21').Value = '  +0.32%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '190.21'
$ws.Range('E22').Value = '  -5.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.966'
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.84'
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('E26').Value = '  -0.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.243'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.95'
$ws.Range('E28').Value = '  -1.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.521'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05643'
$ws.Range('E30').Value = '  -4.11%  '
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.484'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.370'
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.579'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.795'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9462'
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.412'
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5773'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01589'
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.965'
$ws.Range('E40').Value = '  +0.24%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.004'
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8400'
$ws.Range('E42').Value = '  -2.30%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.022.81'
$ws.Range('E43').Value = '  -5.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.17'
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.791.27'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '58.26'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₈104'
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  -1.09%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05331'
$ws.Range('E49').Value = '  +3.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.044'
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4345'
$ws.Range('E51').Value = '  -1.71%  '
